$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Article Identifier(s)" column header is renamed to "Study Identifier".
$ws.Range("A1").Value = "Study Identifier"

# The saved cursor/selection moves from B7 to B1.
$ws.Range("B1").Select()
